# Apply "Doing Updates for Financials" edits to the TNK yearly financials sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TNK")

# --- Income Statement: Earnings Before Interest And Taxes (row 21) ---
$ws.Range("J21").Value = "NA"

# --- Balance Sheet: Long Term Investments (row 47) ---
$ws.Range("D47").Value = 25500
$ws.Range("E47").Value = 70700
$ws.Range("F47").Value = 86800
$ws.Range("G47").Value = 73400

# --- Balance Sheet: Other Assets (row 52) ---
$ws.Range("D52").Value = 7000
$ws.Range("E52").Value = 4600
$ws.Range("F52").Value = 5300
$ws.Range("G52").Value = 8400

# --- Balance Sheet: Accounts Payable (row 57) ---
$ws.Range("D57").Value = 7900
$ws.Range("E57").Value = 14400
$ws.Range("F57").Value = 16700
$ws.Range("G57").Value = 1900
$ws.Range("H57").Value = 2300
$ws.Range("I57").Value = 3300
$ws.Range("J57").Value = 9400

# --- Balance Sheet: Other Current Liabilities (row 59) ---
$ws.Range("D59").Value = 54900
$ws.Range("E59").Value = 97700
$ws.Range("F59").Value = 98900
$ws.Range("G59").Value = 36600
$ws.Range("H59").Value = 42700
$ws.Range("I59").Value = 51100
$ws.Range("J59").Value = 123600

# --- Cash Flow Statement: Depreciation (row 83) ---
$ws.Range("J83").Value = "NA"

# --- Cash Flow Statement: Total Cash Flow From Operating Activities (row 89) ---
$ws.Range("D89").Value = 80500

# --- Cash Flow Statement: Total Cash Flows From Investing Activities (row 94) ---
$ws.Range("J94").Value = "NA"

# --- Cash Flow Statement: Total Cash Flows From Financing Activities (row 100) ---
$ws.Range("D100").Value = -178500
$ws.Range("J100").Value = "NA"

# --- Cash Flow Statement: Effect Of Exchange Rate Changes (row 101) ---
$ws.Range("J101").Value = "NA"

# --- Cash Flow Statement: Change In Cash and Cash Equivalents (row 102) ---
$ws.Range("D102").Value = -19200
